$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3636.2222
$ws.Range("J138").Value = 5252
$ws.Range("L138").Value = 15756
$ws.Range("N138").Value = -26036

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1579.5
$ws.Range("I122").Value = 969.3333
$ws.Range("J122").Value = 2494.75
$ws.Range("K122").Value = 2907.9999
$ws.Range("L122").Value = 7484.25
$ws.Range("M122").Value = -457.9998999999998
$ws.Range("N122").Value = -12384.25

$ws.Range("H123").Value = 29076.334
$ws.Range("I123").Value = 400
$ws.Range("J123").Value = 43414.5
$ws.Range("K123").Value = 400
$ws.Range("L123").Value = 43414.5
$ws.Range("M123").Value = 4500
$ws.Range("N123").Value = -53214.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 43000
$ws.Range("J70").Value = 43000
$ws.Range("L70").Value = 43000
$ws.Range("N70").Value = -43586

$ws.Range("H73").Value = 43000
$ws.Range("J73").Value = 43000
$ws.Range("L73").Value = 43000
$ws.Range("N73").Value = -45028

$ws.Range("H107").Value = 3412.625
$ws.Range("I107").Value = 3471.5715
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 3471.5715
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -1551.5715
$ws.Range("N107").Value = -6840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1992.8182
$ws.Range("I16").Value = 1177.625
$ws.Range("J16").Value = 4166.6665
$ws.Range("K16").Value = 1177.625
$ws.Range("L16").Value = 4166.6665
$ws.Range("M16").Value = -890.625
$ws.Range("N16").Value = -4740.6665

$ws.Range("H58").Value = 2153.96
$ws.Range("I58").Value = 1672.7222
$ws.Range("J58").Value = 3391.4285
$ws.Range("K58").Value = 1672.7222
$ws.Range("L58").Value = 3391.4285
$ws.Range("M58").Value = -1469.7222
$ws.Range("N58").Value = -3797.4285

$ws.Range("H107").Value = 854
$ws.Range("I107").Value = 711.5
$ws.Range("J107").Value = 996.5
$ws.Range("K107").Value = 711.5
$ws.Range("L107").Value = 996.5
$ws.Range("M107").Value = 1208.5
$ws.Range("N107").Value = -4836.5

$ws.Range("H113").Value = 1992.8182
$ws.Range("I113").Value = 1177.625
$ws.Range("J113").Value = 4166.6665
$ws.Range("K113").Value = 1177.625
$ws.Range("L113").Value = 4166.6665
$ws.Range("M113").Value = 992.375
$ws.Range("N113").Value = -8506.666499999999

$ws.Range("H134").Value = 3849.7144
$ws.Range("I134").Value = 1806.1818
$ws.Range("J134").Value = 11342.667
$ws.Range("K134").Value = 5418.5454
$ws.Range("L134").Value = 34028.001
$ws.Range("M134").Value = -2883.5454
$ws.Range("N134").Value = -39098.001

$ws.Range("H136").Value = 2153.96
$ws.Range("I136").Value = 1672.7222
$ws.Range("J136").Value = 3391.4285
$ws.Range("K136").Value = 5018.1666
$ws.Range("L136").Value = 10174.2855
$ws.Range("M136").Value = -2468.1666
$ws.Range("N136").Value = -15274.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3436.5
$ws.Range("I131").Value = 515.7692
$ws.Range("J131").Value = 4955.28
$ws.Range("K131").Value = 1547.3076
$ws.Range("L131").Value = 14865.84
$ws.Range("M131").Value = 3492.6924
$ws.Range("N131").Value = -24945.84

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 25000
$ws.Range("J94").Value = 25000
$ws.Range("L94").Value = 25000
$ws.Range("N94").Value = -26352

$ws.Range("H116").Value = 49644.4
$ws.Range("J116").Value = 49644.4
$ws.Range("L116").Value = 49644.4
$ws.Range("N116").Value = -58822.4

$ws.Range("H122").Value = 11113196
$ws.Range("I122").Value = 14287540
$ws.Range("K122").Value = 42862620
$ws.Range("M122").Value = -42860170

$ws.Range("H132").Value = 7216.2
$ws.Range("I132").Value = 7266.75
$ws.Range("J132").Value = 7014
$ws.Range("K132").Value = 21800.25
$ws.Range("L132").Value = 21042
$ws.Range("M132").Value = -19270.25
$ws.Range("N132").Value = -26102

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 722.9091
$ws.Range("J46").Value = 850.4
$ws.Range("L46").Value = 850.4
$ws.Range("N46").Value = -1226.4

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H88").Value = 9000
$ws.Range("I88").Value = 9000
$ws.Range("K88").Value = 9000
$ws.Range("M88").Value = -8572

$ws.Range("H91").Value = 9000
$ws.Range("I91").Value = 9000
$ws.Range("K91").Value = 9000
$ws.Range("M91").Value = -7518

$ws.Range("H128").Value = 38929
$ws.Range("J128").Value = 38929
$ws.Range("L128").Value = 38929
$ws.Range("N128").Value = -48889

$ws.Range("H136").Value = 3110.2727
$ws.Range("I136").Value = 1568
$ws.Range("J136").Value = 4961
$ws.Range("K136").Value = 4704
$ws.Range("L136").Value = 14883
$ws.Range("M136").Value = -2154
$ws.Range("N136").Value = -19983

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 18500
$ws.Range("J86").Value = 18500
$ws.Range("L86").Value = 18500
$ws.Range("N86").Value = -20746

$ws.Range("H89").Value = 18500
$ws.Range("J89").Value = 18500
$ws.Range("L89").Value = 92500
$ws.Range("N89").Value = -103732

$ws.Range("H107").Value = 496
$ws.Range("I107").Value = 496
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1488
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 432
$ws.Range("N107").ClearContents()

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H125").Value = 35400
$ws.Range("J125").Value = 35400
$ws.Range("L125").Value = 35400
$ws.Range("N125").Value = -45240

$ws.Range("H132").Value = 2566.0386
$ws.Range("I132").Value = 1458.7894
$ws.Range("K132").Value = 4376.3682
$ws.Range("M132").Value = -1846.3682

$ws.Range("H136").Value = 5351.2256
$ws.Range("I136").Value = 5438.1924
$ws.Range("K136").Value = 16314.5772
$ws.Range("M136").Value = -13764.5772
